# Shock.xlsx revert edit
# - Sheet "Y": remove B2/C2 ("Maize (home consumed)" row), fill A3:A6 with 2..5,
#   move selection to A4, sheet no longer the tab-selected one.
# - Sheet "A": becomes the active/selected tab, selection moves to C2
#   (and the view scrolls so row 25 is the top-left row - not reproducible via
#   this COM surface, see notes below).
# - Sheet "VA": selection moves from D10 to D21.
# - Sheet "E": unchanged.

$wb = $excel.ActiveWorkbook

$wsY  = $wb.Worksheets.Item("Y")
$wsA  = $wb.Worksheets.Item("A")
$wsVA = $wb.Worksheets.Item("VA")

# --- Sheet Y: drop the "Maize (home consumed)" row-2 entry, number rows 3-6 ---
$wsY.Range("B2:C2").Clear()

$wsY.Range("A3").Value = 2
$wsY.Range("A4").Value = 3
$wsY.Range("A5").Value = 4
$wsY.Range("A6").Value = 5

# --- Sheet VA: move the remembered selection ---
$wsVA.Select()
$wsVA.Range("D21").Select()

# --- Sheet Y: move its remembered selection (sheet stays inactive afterwards) ---
$wsY.Select()
$wsY.Range("A4").Select()

# --- Sheet A: make it the active tab with its remembered selection ---
$wsA.Select()
$wsA.Range("C2").Select()
